$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08661466666666667
$ws.Range("H2").Value = 0.259844
$ws.Range("I2").Value = 0.5374496355558498
$ws.Range("J2").Value = 0.5374496355558498
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 0.03713543203066667
$ws.Range("R2").Value = 0.334218888276
$ws.Range("S2").Value = 0.002214563345498537
$ws.Range("T2").Value = 0.002214563345498537

$ws.Range("G3").Value = 0.08661466666666667
$ws.Range("H3").Value = 0.259844
$ws.Range("I3").Value = 0.5374496355558498
$ws.Range("J3").Value = 0.5374496355558498
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("Q3").Value = 6.948768775675999
$ws.Range("R3").Value = 62.538918981084
$ws.Range("S3").Value = 0.4143883021005088
$ws.Range("T3").Value = 0.4143883021005088

$ws.Range("G4").Value = 0.08661466666666667
$ws.Range("H4").Value = 0.259844
$ws.Range("I4").Value = 0.5374496355558498
$ws.Range("J4").Value = 0.5374496355558498
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 2.026447799139111
$ws.Range("R4").Value = 18.238030192252
$ws.Range("S4").Value = 0.1208467701098425
$ws.Range("T4").Value = 0.1208467701098425

$ws.Range("I5").Value = 0.2407068810034004
$ws.Range("J5").Value = 0.2407068810034003
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 0.016631798456
$ws.Range("R5").Value = 0.149686186104
$ws.Range("S5").Value = 0.0009918336536373273
$ws.Range("T5").Value = 0.0009918336536373271

$ws.Range("I6").Value = 0.2407068810034004
$ws.Range("J6").Value = 0.2407068810034003
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("Q6").Value = 3.112136185703999
$ws.Range("R6").Value = 28.00922567133599
$ws.Range("S6").Value = 0.185591558955561
$ws.Range("T6").Value = 0.185591558955561

$ws.Range("I7").Value = 0.2407068810034004
$ws.Range("J7").Value = 0.2407068810034003
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 0.9075825844453334
$ws.Range("R7").Value = 8.168243260007999
$ws.Range("S7").Value = 0.05412348839420202
$ws.Range("T7").Value = 0.05412348839420201

$ws.Range("G8").Value = 0.035752
$ws.Range("H8").Value = 0.107256
$ws.Range("I8").Value = 0.2218434834407499
$ws.Range("J8").Value = 0.2218434834407499
$ws.Range("M8").Value = 0.428743
$ws.Range("N8").Value = 1.286229
$ws.Range("O8").Value = 0.00412050394863168
$ws.Range("P8").Value = 0.00412050394863168
$ws.Range("Q8").Value = 0.015328419736
$ws.Range("R8").Value = 0.137955777624
$ws.Range("S8").Value = 0.0009141069494958167
$ws.Range("T8").Value = 0.0009141069494958169

$ws.Range("G9").Value = 0.035752
$ws.Range("H9").Value = 0.107256
$ws.Range("I9").Value = 0.2218434834407499
$ws.Range("J9").Value = 0.2218434834407499
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("O9").Value = 0.7710272268990069
$ws.Range("P9").Value = 0.7710272268990069
$ws.Range("Q9").Value = 2.868248425223999
$ws.Range("R9").Value = 25.814235827016
$ws.Range("S9").Value = 0.1710473658429371
$ws.Range("T9").Value = 0.1710473658429372

$ws.Range("G10").Value = 0.035752
$ws.Range("H10").Value = 0.107256
$ws.Range("I10").Value = 0.2218434834407499
$ws.Range("J10").Value = 0.2218434834407499
$ws.Range("M10").Value = 23.39612766666667
$ws.Range("N10").Value = 70.188383
$ws.Range("O10").Value = 0.2248522691523614
$ws.Range("P10").Value = 0.2248522691523614
$ws.Range("Q10").Value = 0.8364583563386667
$ws.Range("R10").Value = 7.528125207048
$ws.Range("S10").Value = 0.04988201064831693
$ws.Range("T10").Value = 0.04988201064831693
